$d = $word.ActiveDocument

# --- Paragraph 1: "On Pilgrimage - June 1955" (Heading1) -> Title style ---
$p1 = $d.Paragraphs.Item(1)
$p1.Style = "Title"

# --- Paragraph 2: "By Dorothy Day" (bold Normal) -> "Dorothy Day" (Authors style) ---
# Insert a brand-new, unformatted paragraph right after the title so the
# replacement run starts with clean (non-bold) run properties, then delete
# the old "By Dorothy Day" paragraph outright.
$p1.Range.InsertParagraphAfter()
$authors = $d.Paragraphs.Item(2)
$authors.Style = "Authors"
$authors.Range.Text = "Dorothy Day"

$oldAuthors = $d.Paragraphs.Item(3)
$oldAuthors.Range.Delete()
